{"js": "// Replace the 100 arithmetic-problem answers in the single 20x5 table.\n// Each (oldText, newText) pair corresponds to one table cell in row-major\n// order (row 0 col 0..4, row 1 col 0..4, ...). We search *within* each\n// cell's body for its exact current text and do an in-place text replace\n// (InsertLocation.Replace on the found Range) so paragraph/run formatting\n// (font, size, alignment) is left untouched \u2014 only the <w:t> changes.\n\nconst pairs = [[\"76-57=19\", \"68+14=82\"], [\"80-7=73\", \"30-6=24\"], [\"29+46=75\", \"14+39=53\"], [\"15+66=81\", \"90-69=21\"], [\"47+6=53\", \"58+5=63\"], [\"50-42=8\", \"59+32=91\"], [\"18+75=93\", \"71-13=58\"], [\"79+2=81\", \"88-29=59\"], [\"39+43=82\", \"85-78=7\"], [\"93-39=54\", \"75-47=28\"], [\"66-27=39\", \"76-38=38\"], [\"87+5=92\", \"83-78=5\"], [\"5+76=81\", \"5+8=13\"], [\"56-28=28\", \"87-38=49\"], [\"61-8=53\", \"49+5=54\"], [\"79+2=81\", \"27+14=41\"], [\"67+9=76\", \"8+69=77\"], [\"73+18=91\", \"27+24=51\"], [\"9+17=26\", \"95-39=56\"], [\"48+49=97\", \"65-46=19\"], [\"2+9=11\", \"9+82=91\"], [\"49+18=67\", \"29+16=45\"], [\"19+4=23\", \"61-27=34\"], [\"35-6=29\", \"55+37=92\"], [\"56-47=9\", \"55-16=39\"], [\"55+19=74\", \"9+64=73\"], [\"51-4=47\", \"47+17=64\"], [\"80-69=11\", \"2+79=81\"], [\"40-39=1\", \"6+7=13\"], [\"71-56=15\", \"64+27=91\"], [\"34+9=43\", \"65+9=74\"], [\"80-11=69\", \"74-69=5\"], [\"16+56=72\", \"39+19=58\"], [\"37+55=92\", \"54-26=28\"], [\"28+48=76\", \"94-85=9\"], [\"16+39=55\", \"9+52=61\"], [\"17-8=9\", \"6+19=25\"], [\"21-12=9\", \"21-13=8\"], [\"8+89=97\", \"93-76=17\"], [\"48+29=77\", \"19+39=58\"], [\"59+6=65\", \"37-8=29\"], [\"53-49=4\", \"69+23=92\"], [\"29+7=36\", \"89+2=91\"], [\"45+47=92\", \"7+16=23\"], [\"37+7=44\", \"56+18=74\"], [\"8+75=83\", \"90-12=78\"], [\"7+77=84\", \"83-49=34\"], [\"26+28=54\", \"17+46=63\"], [\"59+38=97\", \"45-18=27\"], [\"83-15=68\", \"97-59=38\"], [\"95-16=79\", \"94-35=59\"], [\"63-14=49\", \"84+7=91\"], [\"70-5=65\", \"90-65=25\"], [\"33+28=61\", \"90-66=24\"], [\"33+49=82\", \"39+26=65\"], [\"42-5=37\", \"63-15=48\"], [\"47+18=65\", \"80-71=9\"], [\"59+16=75\", \"94-65=29\"], [\"60-48=12\", \"26+49=75\"], [\"69+25=94\", \"62-14=48\"], [\"78+16=94\", \"71-49=22\"], [\"90-25=65\", \"9+49=58\"], [\"19+17=36\", \"49+36=85\"], [\"9+46=55\", \"27+24=51\"], [\"27+16=43\", \"80-67=13\"], [\"52-16=36\", \"71-35=36\"], [\"17+28=45\", \"8+47=55\"], [\"87-19=68\", \"58+3=61\"], [\"48+16=64\", \"72-64=8\"], [\"92-86=6\", \"90-78=12\"], [\"63+29=92\", \"66+7=73\"], [\"64-35=29\", \"16+37=53\"], [\"12+39=51\", \"35-26=9\"], [\"63-26=37\", \"19+67=86\"], [\"92-79=13\", \"62-55=7\"], [\"76-48=28\", \"51-9=42\"], [\"30-16=14\", \"49+29=78\"], [\"75+8=83\", \"16+79=95\"], [\"21-17=4\", \"6+59=65\"], [\"53-45=8\", \"48+4=52\"], [\"91-69=22\", \"19+15=34\"], [\"60-45=15\", \"33-28=5\"], [\"71-52=19\", \"81-27=54\"], [\"72-14=58\", \"27+4=31\"], [\"49+37=86\", \"48+5=53\"], [\"63-58=5\", \"51-42=9\"], [\"40-38=2\", \"35+28=63\"], [\"80-25=55\", \"55-49=6\"], [\"80-32=48\", \"12-9=3\"], [\"54-16=38\", \"47+7=54\"], [\"73-14=59\", \"84-28=56\"], [\"40-16=24\", \"72-28=44\"], [\"60-25=35\", \"78+5=83\"], [\"41-14=27\", \"74-49=25\"], [\"6+5=11\", \"60-27=33\"], [\"9+48=57\", \"90-81=9\"], [\"35+8=43\", \"5+39=44\"], [\"51-44=7\", \"6+18=24\"], [\"48+18=66\", \"47-29=18\"], [\"80-36=44\", \"35+9=44\"]];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst columnCount = table.values[0].length;\n\nlet index = 0;\nfor (let row = 0; row < table.rowCount; row++) {\n  for (let col = 0; col < columnCount; col++) {\n    const [oldText, newText] = pairs[index];\n    const cell = table.getCell(row, col);\n    const results = cell.body.search(oldText, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length > 0) {\n      results.items[0].insertText(newText, Word.InsertLocation.replace);\n    }\n    index++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem answers in the single 20x5 table.\n# Each pair is (oldText, newText) for one cell, in row-major order\n# (row 1 col 1..5, row 2 col 1..5, ...). Setting Cell.Range.Text only\n# rewrites the run's <w:t> content; paragraph/run formatting (font,\n# size, alignment) defined in pPr/rPr is left untouched.\n\n$pairs = @(\n    ,@(\"76-57=19\", \"68+14=82\")\n    ,@(\"80-7=73\", \"30-6=24\")\n    ,@(\"29+46=75\", \"14+39=53\")\n    ,@(\"15+66=81\", \"90-69=21\")\n    ,@(\"47+6=53\", \"58+5=63\")\n    ,@(\"50-42=8\", \"59+32=91\")\n    ,@(\"18+75=93\", \"71-13=58\")\n    ,@(\"79+2=81\", \"88-29=59\")\n    ,@(\"39+43=82\", \"85-78=7\")\n    ,@(\"93-39=54\", \"75-47=28\")\n    ,@(\"66-27=39\", \"76-38=38\")\n    ,@(\"87+5=92\", \"83-78=5\")\n    ,@(\"5+76=81\", \"5+8=13\")\n    ,@(\"56-28=28\", \"87-38=49\")\n    ,@(\"61-8=53\", \"49+5=54\")\n    ,@(\"79+2=81\", \"27+14=41\")\n    ,@(\"67+9=76\", \"8+69=77\")\n    ,@(\"73+18=91\", \"27+24=51\")\n    ,@(\"9+17=26\", \"95-39=56\")\n    ,@(\"48+49=97\", \"65-46=19\")\n    ,@(\"2+9=11\", \"9+82=91\")\n    ,@(\"49+18=67\", \"29+16=45\")\n    ,@(\"19+4=23\", \"61-27=34\")\n    ,@(\"35-6=29\", \"55+37=92\")\n    ,@(\"56-47=9\", \"55-16=39\")\n    ,@(\"55+19=74\", \"9+64=73\")\n    ,@(\"51-4=47\", \"47+17=64\")\n    ,@(\"80-69=11\", \"2+79=81\")\n    ,@(\"40-39=1\", \"6+7=13\")\n    ,@(\"71-56=15\", \"64+27=91\")\n    ,@(\"34+9=43\", \"65+9=74\")\n    ,@(\"80-11=69\", \"74-69=5\")\n    ,@(\"16+56=72\", \"39+19=58\")\n    ,@(\"37+55=92\", \"54-26=28\")\n    ,@(\"28+48=76\", \"94-85=9\")\n    ,@(\"16+39=55\", \"9+52=61\")\n    ,@(\"17-8=9\", \"6+19=25\")\n    ,@(\"21-12=9\", \"21-13=8\")\n    ,@(\"8+89=97\", \"93-76=17\")\n    ,@(\"48+29=77\", \"19+39=58\")\n    ,@(\"59+6=65\", \"37-8=29\")\n    ,@(\"53-49=4\", \"69+23=92\")\n    ,@(\"29+7=36\", \"89+2=91\")\n    ,@(\"45+47=92\", \"7+16=23\")\n    ,@(\"37+7=44\", \"56+18=74\")\n    ,@(\"8+75=83\", \"90-12=78\")\n    ,@(\"7+77=84\", \"83-49=34\")\n    ,@(\"26+28=54\", \"17+46=63\")\n    ,@(\"59+38=97\", \"45-18=27\")\n    ,@(\"83-15=68\", \"97-59=38\")\n    ,@(\"95-16=79\", \"94-35=59\")\n    ,@(\"63-14=49\", \"84+7=91\")\n    ,@(\"70-5=65\", \"90-65=25\")\n    ,@(\"33+28=61\", \"90-66=24\")\n    ,@(\"33+49=82\", \"39+26=65\")\n    ,@(\"42-5=37\", \"63-15=48\")\n    ,@(\"47+18=65\", \"80-71=9\")\n    ,@(\"59+16=75\", \"94-65=29\")\n    ,@(\"60-48=12\", \"26+49=75\")\n    ,@(\"69+25=94\", \"62-14=48\")\n    ,@(\"78+16=94\", \"71-49=22\")\n    ,@(\"90-25=65\", \"9+49=58\")\n    ,@(\"19+17=36\", \"49+36=85\")\n    ,@(\"9+46=55\", \"27+24=51\")\n    ,@(\"27+16=43\", \"80-67=13\")\n    ,@(\"52-16=36\", \"71-35=36\")\n    ,@(\"17+28=45\", \"8+47=55\")\n    ,@(\"87-19=68\", \"58+3=61\")\n    ,@(\"48+16=64\", \"72-64=8\")\n    ,@(\"92-86=6\", \"90-78=12\")\n    ,@(\"63+29=92\", \"66+7=73\")\n    ,@(\"64-35=29\", \"16+37=53\")\n    ,@(\"12+39=51\", \"35-26=9\")\n    ,@(\"63-26=37\", \"19+67=86\")\n    ,@(\"92-79=13\", \"62-55=7\")\n    ,@(\"76-48=28\", \"51-9=42\")\n    ,@(\"30-16=14\", \"49+29=78\")\n    ,@(\"75+8=83\", \"16+79=95\")\n    ,@(\"21-17=4\", \"6+59=65\")\n    ,@(\"53-45=8\", \"48+4=52\")\n    ,@(\"91-69=22\", \"19+15=34\")\n    ,@(\"60-45=15\", \"33-28=5\")\n    ,@(\"71-52=19\", \"81-27=54\")\n    ,@(\"72-14=58\", \"27+4=31\")\n    ,@(\"49+37=86\", \"48+5=53\")\n    ,@(\"63-58=5\", \"51-42=9\")\n    ,@(\"40-38=2\", \"35+28=63\")\n    ,@(\"80-25=55\", \"55-49=6\")\n    ,@(\"80-32=48\", \"12-9=3\")\n    ,@(\"54-16=38\", \"47+7=54\")\n    ,@(\"73-14=59\", \"84-28=56\")\n    ,@(\"40-16=24\", \"72-28=44\")\n    ,@(\"60-25=35\", \"78+5=83\")\n    ,@(\"41-14=27\", \"74-49=25\")\n    ,@(\"6+5=11\", \"60-27=33\")\n    ,@(\"9+48=57\", \"90-81=9\")\n    ,@(\"35+8=43\", \"5+39=44\")\n    ,@(\"51-44=7\", \"6+18=24\")\n    ,@(\"48+18=66\", \"47-29=18\")\n    ,@(\"80-36=44\", \"35+9=44\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$columnCount = $t.Columns.Count\n\n$index = 0\nfor ($row = 1; $row -le $t.Rows.Count; $row++) {\n    for ($col = 1; $col -le $columnCount; $col++) {\n        $newText = $pairs[$index][1]\n        $t.Cell($row, $col).Range.Text = $newText\n        $index++\n    }\n}\n\nWrite-Output \"updated $index cells\"\n"}
